$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.115.66'
$ws.Range('D3').Value = '2.315.69'
$ws.Range('E3').Value = '  -2.24%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '313.26'
$ws.Range('E5').Value = '  -5.87%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '106.46'
$ws.Range('E6').Value = '  +4.76%  '
$ws.Range('E7').Value = '  -2.26%  '
$ws.Range('E9').Value = '  -2.87%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '40.57'
$ws.Range('E10').Value = '  +1.22%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0917'
$ws.Range('E11').Value = '  -0.73%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '8.31'
$ws.Range('E12').Value = '  -1.92%  '
$ws.Range('E13').Value = '  +0.30%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.985'
$ws.Range('E14').Value = '  -2.21%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '15.67'
$ws.Range('E15').Value = '  -4.89%  '
$ws.Range('D16').Value = '2.668.73'
$ws.Range('E16').Value = '  -2.07%  '
$ws.Range('D17').Value = '2.352.65'
$ws.Range('E17').Value = '  -0.62%  '
$ws.Range('D18').Value = '42.164.41'
$ws.Range('E18').Value = '  -1.27%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.69'
$ws.Range('E19').Value = '  -5.25%  '
$ws.Range('E20').Value = '  -1.58%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '74.85'
$ws.Range('E21').Value = '  -2.61%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '3.50'
$ws.Range('E22').Value = '  -7.02%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '261.05'
$ws.Range('E23').Value = '  -3.15%  '
$ws.Range('E24').Value = '  -0.65%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.33'
$ws.Range('E25').Value = '  -8.40%  '
$ws.Range('E26').Value = '  +0.60%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '11.03'
$ws.Range('E27').Value = '  -4.31%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '22.86'
$ws.Range('E28').Value = '  -1.37%  '
$ws.Range('E29').Value = '  +1.33%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '35.91'
$ws.Range('E30').Value = '  +1.53%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0902'
$ws.Range('E31').Value = '  +0.01%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '163.87'
$ws.Range('E32').Value = '  -7.10%  '
$ws.Range('E33').Value = '  -6.14%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.89'
$ws.Range('E34').Value = '  -3.80%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.120'
$ws.Range('E35').Value = '  +13.22%  '
$ws.Range('B36').Value = 'Stellar'
$ws.Range('C36').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.130'
$ws.Range('E36').Value = '  -2.00%  '
$ws.Range('E37').Value = '  -1.81%  '
$ws.Range('E38').Value = '  -1.36%  '
$ws.Range('E39').Value = '  -6.39%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.64'
$ws.Range('E40').Value = '  -4.37%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '98.75'
$ws.Range('E41').Value = '  +6.55%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '71.41'
$ws.Range('E42').Value = '  +1.81%  '
$ws.Range('E43').Value = '  -3.57%  '
$ws.Range('E44').Value = '  -2.25%  '
$ws.Range('E45').Value = '  -0.11%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '12.46'
$ws.Range('E46').Value = '  +5.52%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '112.25'
$ws.Range('E47').Value = '  -5.16%  '
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '9.03'
$ws.Range('E48').Value = '  -2.02%  '
$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '5.37'
$ws.Range('E49').Value = '  -2.58%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '74.65'
$ws.Range('E50').Value = '  +4.09%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.27'
$ws.Range('E51').Value = '  -0.34%  '
